# This script applies a cyclic rotation of several paragraph contents in the
# document, matching the commit's diff:
#
#   Objetivos text          <- old "Programa resumido" text (short variant)
#   Docente text             <- old "Objetivos" text
#   Programa resumido text   <- old "Programa" text (long variant, typo "Geotéxteis")
#   Programa text            <- old "Método" value
#   Método value              <- old "Critério" value
#   Critério value            <- old "Norma de recuperação" value
#   Norma de recuperação      <- old "Bibliografia" full content (multi-line, w:br separated)
#   Bibliografia content      <- old "Docente" text
#
# Because this is a closed rotation cycle, every original value is first
# swapped out for a unique placeholder token; only afterwards are the
# placeholders swapped in for their final values. This guarantees no
# replacement step can accidentally re-match text that a previous step just
# inserted.

$d = $word.ActiveDocument

function Replace-Literal($findText, $replaceText) {
    $ok = $d.Content.Find.Execute(
        $findText,    # FindText
        $true,        # MatchCase
        $true,        # MatchWholeWord
        $false,       # MatchWildcards
        $false,       # MatchSoundsLike
        $false,       # MatchAllWordForms
        $true,        # Forward
        1,            # Wrap (wdFindContinue)
        $false,       # Format
        $replaceText, # ReplaceWith
        2             # Replace (wdReplaceAll)
    )
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# ---- Original (pre-edit) paragraph values ----
$oldObjetivos = "Apresentar ao aluno as aplicações dos geossintéticos em obras de proteção e recuperação ambiental dando ênfase às múltiplas funções dos geossintéticos. Aspectos como a drenagem e filtração, a separação, barreiras de silte, reforço de solos e os sistemas de contenção de resíduos e de efluentes que podem ser utilizados com vantagens técnicas e econômicas tanto em obras de proteção ambiental como na recuperação de áreas degradadas e/ou contaminadas."
$oldDocente = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$oldResumido = "Obras de proteção ambiental: sistemas de contenção de resíduos; lagoas de efluentes, remediação e mitigação de áreas degradadas; Geossintéticos: tipos e aplicações; Geotêxteis, geomantas e geocompostos para a drenagem: tipos, propriedades, ensaios caracterização e desempenho; Drenagem e filtração com geossintéticos. Critérios de filtração com geossintéticos. Detalhes construtivos de obras de drenagem com geossintéticos; Sistemas de drenagem em obras de proteção e de recuperação ambiental; Geogrelhas e geotêxteis para reforço: tipos, propriedades, ensaios de caracterização e de desempenho; Princípios gerais de reforço de solos. Reforço de solos com geossintéticos. Estruturas em solos reforçados. Detalhamento da construção de obras em solo reforçado; Geomembranas: tipos, propriedades, ensaios de caracterização e de desempenho; Geocomposto bentonítico: tipos, propriedades, ensaios de caracterização e de desempenho; Barreiras impermeabilizantes de fundação e de cobertura."
$oldPrograma = "Obras de proteção ambiental: sistemas de contenção de resíduos; lagoas de efluentes, remediação e mitigação de áreas degradadas; Geossintéticos: tipos e aplicações; Geotéxteis, geomantas e geocompostos para a drenagem: tipos, propriedades; ensaios de caracterização e de desempenho; Drenagem e filtração com geossintéticos. Critérios de filtração com geossintéticos. Detalhes construtivos de obras de drenagem com geossintéticos; Sistemas de drenagem em obras de proteção e de recuperação ambiental; Geogrelhas e geotêxteis para reforço; tipos, propriedades, ensaios de caracterização e de desempenho; Princípios gerais de reforço de solos. Reforço de solos com geossintéticos. Estruturas em solos reforçados. Detalhamento da construção de obras em solo reforçado; Geomembranas: tipos, propriedades, ensaios de caracterização e de desempenho; Geocomposto bentonítico: tipos, propriedades, ensaios de caracterização e de desempenho; Barreiras impermeabilizantes de fundação e de cobertura: tipos, especificações, características; Solicitações físicas mecânicas e químicas de geomembranas em sistemas de contenção de resíduos e de efluentes; Aspectos construtivos de barreiras impermeabilizantes e de cobertura; especificações de geossintéticos."
$oldMetodo = "Aulas expositivas; microcomputadores; seminários; visitas técnicas."
$oldCriterio = "Média ponderada de notas de provas e seminários."
$oldNorma = "Prova única com nota igual ou superior a 5,0 (cinco)."

# The original Bibliografia paragraph spans a single run broken up by
# <w:br/> manual line breaks. In Word's Range.Text / Find story text, each
# <w:br/> surfaces as a vertical-tab character (`v in PowerShell escaping),
# so the whole paragraph can be matched (and replaced) in one Find.Execute.
$oldBibliografiaFull = "SHARMA, H. D. & LEWIS, S.P. (1994) Waste Containement Systems, Waste Stabilization and Landfils. Design and Evaluation - Joh Willy & Sons, Inc. New York;`vKOERNER. R.M. (1997) Designing with Geosynthetics, Prentice Hall Inc.;`vQIAN, X/; KOERNER, R.M. & GRAY, D.H. (2002) Geotechnical Aspects of Landfill Design and Construction - Prentice Hall, Upper Saddle Riner, New Jersey 07458;`vManual Brasileiro de Geossintéticos, ABINT, 2004;`vINGOLD, T.S. (1994) The Geotêxtiles and Geomembranes. Manual, Elsevier, London, 610p.;`vVAN ZATEN, R.V. (1986) Geotextiles and Geomembranes in Civil Engineering. A.A. Balkema, Rotherdam, netherlands, 654p."

# Replacement text for the new "Norma de recuperação" value: uses "^l" (Word
# Find/Replace group syntax for a manual line break) so the resulting run
# ends up with the same <w:t>/<w:br/> structure as the original Bibliografia
# paragraph.
$bibliografiaAsNormaValue = "SHARMA, H. D. & LEWIS, S.P. (1994) Waste Containement Systems, Waste Stabilization and Landfils. Design and Evaluation - Joh Willy & Sons, Inc. New York;^lKOERNER. R.M. (1997) Designing with Geosynthetics, Prentice Hall Inc.;^lQIAN, X/; KOERNER, R.M. & GRAY, D.H. (2002) Geotechnical Aspects of Landfill Design and Construction - Prentice Hall, Upper Saddle Riner, New Jersey 07458;^lManual Brasileiro de Geossintéticos, ABINT, 2004;^lINGOLD, T.S. (1994) The Geotêxtiles and Geomembranes. Manual, Elsevier, London, 610p.;^lVAN ZATEN, R.V. (1986) Geotextiles and Geomembranes in Civil Engineering. A.A. Balkema, Rotherdam, netherlands, 654p."

# ---- Placeholder tokens (unique, will not collide with any real content) ----
$P_OBJETIVOS = "@@PLACEHOLDER_OBJETIVOS@@"
$P_DOCENTE   = "@@PLACEHOLDER_DOCENTE@@"
$P_RESUMIDO  = "@@PLACEHOLDER_RESUMIDO@@"
$P_PROGRAMA  = "@@PLACEHOLDER_PROGRAMA@@"
$P_METODO    = "@@PLACEHOLDER_METODO@@"
$P_CRITERIO  = "@@PLACEHOLDER_CRITERIO@@"
$P_NORMA     = "@@PLACEHOLDER_NORMA@@"
$P_BIBL      = "@@PLACEHOLDER_BIBLIOGRAFIA@@"

# ---- Phase 1: swap every original value out for its own placeholder ----
Replace-Literal $oldObjetivos        $P_OBJETIVOS
Replace-Literal $oldDocente          $P_DOCENTE
Replace-Literal $oldResumido         $P_RESUMIDO
Replace-Literal $oldPrograma         $P_PROGRAMA
Replace-Literal $oldMetodo           $P_METODO
Replace-Literal $oldCriterio         $P_CRITERIO
Replace-Literal $oldNorma            $P_NORMA
Replace-Literal $oldBibliografiaFull $P_BIBL

# ---- Phase 2: swap placeholders in for their final (rotated) values ----
Replace-Literal $P_OBJETIVOS $oldResumido
Replace-Literal $P_DOCENTE   $oldObjetivos
Replace-Literal $P_RESUMIDO  $oldPrograma
Replace-Literal $P_PROGRAMA  $oldMetodo
Replace-Literal $P_METODO    $oldCriterio
Replace-Literal $P_CRITERIO  $oldNorma
Replace-Literal $P_NORMA     $bibliografiaAsNormaValue
Replace-Literal $P_BIBL      $oldDocente

Write-Output "Done."
